$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("translations")

$ws.Range("A42").Value = "watchlist_only"
$ws.Range("B42").Value = "nur Einträge in der Merkliste anzeigen"

$ws.Range("A43").Value = "records"
$ws.Range("B43").Value = "Einträge"

$ws.Range("A44").Value = "record"
$ws.Range("B44").Value = "Eintrag"
